# Update reviewdb worksheet (Sheet1):
#  - row 40: "confirm" -> "no" in column G
#  - row 41: shift review-chain values (C/F/G) by one position, replacing the
#    removed "confirm" entry, so C/F/G pick up the next email/review/flag
#  - rows 42 and 43: two new review rows are appended (previously blank
#    placeholder rows), each with its own appid/keyword/email/recovery/time/
#    review/flag values
#  - selection moves from G27 to G44

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40: flag column updated ---
$ws.Range("G40").Value = "no"

# --- Row 41: email/review/flag columns updated ---
$ws.Range("C41").Value = "dony1098765432@gmail.com"
$ws.Range("F41").Value = "welcome to the best guide about bitcoin this year"
$ws.Range("G41").Value = "no"

# --- Rows 42 & 43 are new data rows; copy formatting from row 41 first ---
$ws.Range("A41:G41").Copy()
$ws.Range("A42:G42").PasteSpecial(-4122)
$ws.Range("A41:G41").Copy()
$ws.Range("A43:G43").PasteSpecial(-4122)

$ws.Range("A42").Value = "com.hamxa.shaynachim"
$ws.Range("B42").Value = "bitcoin"
$ws.Range("C42").Value = "cohenn167@gmail.com"
$ws.Range("D42").Value = "stavsade45@gmail.com"
$ws.Range("E42").Value = "27/5/2019 15:59"
$ws.Range("F42").Value = "nice car and tracks! Like it"
$ws.Range("G42").Value = "no"

$ws.Range("A43").Value = "com.hamxa.shaynachim"
$ws.Range("B43").Value = "bitcoin"
$ws.Range("C43").Value = "cohenyossi408@gmail.com"
$ws.Range("D43").Value = "cohenn167@gmail.com"
$ws.Range("E43").Value = "27/5/2019 15:59"
$ws.Range("F43").Value = "awesome app with great addictive concept"
$ws.Range("G43").Value = "no"

# --- Update the active selection to G44 ---
$ws.Range("G44").Select()
